$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Important Notes" text in A1: middle school PD hours 65 -> 55
$notes = $ws.Range("A1").Value2
$notes = $notes -replace "total hours for middle school PD is 65", "total hours for middle school PD is 55"
$ws.Range("A1").Value = $notes

# Update the formulas that compute Middle School costs: 65 -> 55
$ws.Range("C8").Formula = "=(55*B3*B5)/2"
$ws.Range("C9").Formula = "=(55*B3*B5)/2"

# Move the selection to C10
$ws.Range("C10").Select()
